$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "sex" coded as 2 = female, 1 = male.
# Recode the numeric codes to their text labels.
# NOTE: female rows are written first so the new shared-string table
# picks up "female" before "male" (matches the source ordering).
$femaleRows = @(2,8,9,10,13,14,18,20,22,27,32,33,37,41,46,47,49,51,52,53,60,61,64,67,72,74,77,80,82,85,86,91,92,94,96,99,101,102,104,106,107,109,115,117,119,121,124,127,128,131,134,138,139,140,141,143,145,146,147,148,151)
$maleRows = @(3,4,5,6,7,11,12,15,16,17,19,21,23,24,25,26,28,29,30,31,34,35,36,38,39,40,42,43,44,45,48,50,54,55,56,57,58,59,62,63,65,66,68,69,70,71,73,75,76,78,79,81,83,84,87,88,89,90,93,95,97,98,100,103,105,108,110,111,112,113,114,116,118,120,122,123,125,126,129,130,132,133,135,136,137,142,144,149,150,152)

foreach ($r in $femaleRows) {
    $ws.Cells.Item($r, 4).Value = "female"
}
foreach ($r in $maleRows) {
    $ws.Cells.Item($r, 4).Value = "male"
}

# A couple of unrelated data-entry fixes further down the sheet.
$ws.Cells.Item(10, 6).Value = 108
$ws.Cells.Item(29, 12).Value = 3.67

# Restore the view: scroll back to A1 and leave the active selection on Q25.
$ws.Range("Q25").Select()
